$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("web 120")

$ws.Range("D22").Value = 58
$ws.Range("D23").Value = 60
$ws.Range("D27").Value = 50
$ws.Range("D29").Value = 60

$ws.Range("H26").Select()

$ws.Columns("F").AutoFit() | Out-Null

$wb.Save()
